$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix corrupted date formats in column B (rows 6-11): MM/DD/YYYY -> YYYY-MM-DD ---
# Force text formatting first so Excel doesn't reinterpret the new value as a date,
# then clear the formatting again so the cell keeps its original (unstyled) look.
$dateRange = $ws.Range("B6:B11")
$dateRange.NumberFormat = "@"

$ws.Range("B6").Value = "2022-10-23"
$ws.Range("B7").Value = "2022-11-28"
$ws.Range("B8").Value = "2022-11-28"
$ws.Range("B9").Value = "2022-12-10"
$ws.Range("B10").Value = "2022-12-10"
$ws.Range("B11").Value = "2022-12-10"

$dateRange.ClearFormats()

# --- Append new row 12 ("methods" dataset) ---
$ws.Range("A12").Value = "0f6bad26-cd71-4c67-ae5e-66bd7c37fc02"
# Reuse the same formatting as the other ID cells in column A (bordered/bold/centered style).
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "2022-12-10"
$ws.Range("B12").ClearFormats()

# C12 stays blank (matches the empty "modified" cells above it) - a leading
# apostrophe keeps it a real (empty) text cell instead of clearing it outright.
$ws.Range("C12").Value = "'"
$ws.Range("C12").ClearFormats()

$ws.Range("D12").Value = "methods"
$ws.Range("E12").Value = "methods.xlsx"
$ws.Range("F12").Value = "active"
